# Adição da Sprint 9 e Alterações no Burndown
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab (capitalize "Sprint 1")
$ws.Name = "Sprint 1"

# Hide gridlines on the sheet view
$excel.ActiveWindow.DisplayGridlines = $false

# Insert a new "HORAS TRABALHADAS" column before STATUS (currently column D)
# so STATUS shifts right to column E.
$ws.Range("D1").EntireColumn.Insert()

$ws.Range("D1").Value2 = "HORAS TRABALHADAS"
$ws.Range("D2").Value2 = " 3 Horas e 15 Minutos"
$ws.Range("D3").Value2 = " 3 Horas e 15 Minutos"

# Match the horizontal-center-only alignment already used by column B
# (nudge VerticalAlignment off its inherited "center" so the engine
# actually swaps in the center/bottom style instead of treating it as a no-op)
$ws.Range("D2:D3").VerticalAlignment = -4107

# Widen the new column like the rest of the table
$ws.Columns.Item(4).ColumnWidth = 19.6

# Rebuild the table (ListObject) over the new range so the column list,
# header names and table style all resync correctly.
$tbl = $ws.ListObjects.Item(1)
$tbl.Unlist()
$tbl2 = $ws.ListObjects.Add(1, $ws.Range("A1:E3"), [System.Reflection.Missing]::Value, 1)
$tbl2.Name = "Tabela1"
$tbl2.TableStyle = "TableStyleMedium6"

# Match the final selection/active cell from the authored workbook
$ws.Range("D8").Select()
